# StagingTemplates/Staging.OutputPersonRole.xlsx
#
# The workbook's "OutputPersonRole" staging header row is being renamed
# from the old *SourceKey naming convention to the new *BusinessKey
# convention, and a fourth header ("OutputBusinessKey") is introduced
# in column A (shifting the previous column-A header into column B).
#
# Resulting header row (row 2):
#   A2 = OutputBusinessKey
#   B2 = OutputPersonRole_ID
#   C2 = PersonBusinessKey
#   D2 = RoleBusinessKey

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "OutputBusinessKey"
$ws.Range("B2").Value = "OutputPersonRole_ID"
$ws.Range("C2").Value = "PersonBusinessKey"
$ws.Range("D2").Value = "RoleBusinessKey"
